$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, centered, bordered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-4: new columns I (value 1) and J (value 4)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4
